$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted like "1.234.56" (thousand-sep dots) which
# Excel would otherwise auto-convert to a number when assigned via .Value. Force
# these specific cells to Text format first so the literal string is preserved,
# exactly like the rest of the (already-text) column.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D12', 'D13', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D31', 'D32', 'D35', 'D39', 'D40', 'D42', 'D43', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.666.20'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').Value = '3.369.42'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '567.09'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').Value = '135.61'
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.366.97'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').Value = '3.942.98'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '25.51'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').Value = '3.368.71'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '60.902.87'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '5.77'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '13.68'
$ws.Range('E20').Value = '  -3.44%  '
$ws.Range('D21').Value = '9.26'
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('D22').Value = '372.12'
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').Value = '3.510.98'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').Value = '0.544'
$ws.Range('E24').Value = '  -2.73%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '70.72'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('E28').Value = '  +11.16%  '
$ws.Range('E29').Value = '  -6.13%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').Value = '7.32'
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  -1.92%  '
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '23.18'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('E36').Value = '  -4.35%  '
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').Value = '164.67'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').Value = '0.0753'
$ws.Range('E40').Value = '  -3.92%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '0.771'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D43').Value = '24.91'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('E46').Value = '  -5.74%  '
$ws.Range('D47').Value = '2.527.61'
$ws.Range('E47').Value = '  +8.36%  '
$ws.Range('D48').Value = '23.00'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '6.74'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('E50').Value = '  +4.24%  '
$ws.Range('D51').Value = '0.0256'
$ws.Range('E51').Value = '  -1.74%  '
